$wb = $excel.ActiveWorkbook

# The status "Ready for handoff" is a shared string reused by the Status
# column (C) on all three sheets for this row (record
# 363e1d64-5b85-4359-b4c1-0676f87ff528). It now becomes "Handback transform
# failed" everywhere it is shown, so update every occurrence identically
# (Overview also duplicates Status into column B) so they collapse back into
# a single shared string on save.
$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# Sheet "zh-cn": Row 3 status + new cell L3 with handback-mismatch error detail
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("L3").Value = "Handback file name: mo55quvd.gky is different with handoff file name: 363e1d64-5b85-4359-b4c1-0676f87ff528.1fddf20ee25bd9a022ec7a961af27d9c56c0e598.zh-cn."

# Sheet "de-de": Row 3 status + new cell L3 with handback-mismatch error detail
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("L3").Value = "Handback file name: mo55quvd.gky is different with handoff file name: 363e1d64-5b85-4359-b4c1-0676f87ff528.1fddf20ee25bd9a022ec7a961af27d9c56c0e598.de-de."
